$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update District (column G) values to official names
$ws.Range("G3").Value = 'Ballari (Bellary)'
$ws.Range("G5").Value = 'Mysuru (Mysore)'
$ws.Range("G6").Value = 'Mysuru (Mysore)'
$ws.Range("G7").Value = 'Kalaburagi (Gulbarga)'
$ws.Range("G8").Value = 'Ballari (Bellary)'
$ws.Range("G9").Value = 'Mysuru (Mysore)'
$ws.Range("G10").Value = 'Mysuru (Mysore)'
$ws.Range("G11").Value = 'Mysuru (Mysore)'
$ws.Range("G12").Value = 'Kalaburagi (Gulbarga)'
$ws.Range("G13").Value = 'Mysuru (Mysore)'
$ws.Range("G14").Value = 'Ballari (Bellary)'
$ws.Range("G15").Value = 'Mysuru (Mysore)'
$ws.Range("G16").Value = 'Mysuru (Mysore)'
$ws.Range("G17").Value = 'Mysuru (Mysore)'
$ws.Range("G19").Value = 'Mysuru (Mysore)'
$ws.Range("G20").Value = 'Ballari (Bellary)'
$ws.Range("G21").Value = 'Kalaburagi (Gulbarga)'
$ws.Range("G22").Value = 'Ballari (Bellary)'
$ws.Range("G23").Value = 'Ballari (Bellary)'
$ws.Range("G24").Value = 'Mysuru (Mysore)'
$ws.Range("G25").Value = 'Mysuru (Mysore)'
$ws.Range("G26").Value = 'Kalaburagi (Gulbarga)'
$ws.Range("G27").Value = 'Mysuru (Mysore)'
$ws.Range("G28").Value = 'Mysuru (Mysore)'
$ws.Range("G29").Value = 'Mysuru (Mysore)'
$ws.Range("G30").Value = 'Mysuru (Mysore)'
$ws.Range("G31").Value = 'Mysuru (Mysore)'
$ws.Range("G32").Value = 'Kalaburagi (Gulbarga)'
$ws.Range("G33").Value = 'Ballari (Bellary)'
$ws.Range("G34").Value = 'Ballari (Bellary)'
$ws.Range("G36").Value = 'Ballari (Bellary)'
$ws.Range("G37").Value = 'Ballari (Bellary)'
$ws.Range("G39").Value = 'Mysuru (Mysore)'
$ws.Range("G40").Value = 'Ballari (Bellary)'
$ws.Range("G41").Value = 'Mysuru (Mysore)'
$ws.Range("G42").Value = 'Mysuru (Mysore)'
$ws.Range("G43").Value = 'Ballari (Bellary)'
$ws.Range("G44").Value = 'Mysuru (Mysore)'
$ws.Range("G45").Value = 'Mysuru (Mysore)'
$ws.Range("G46").Value = 'Ballari (Bellary)'
$ws.Range("G47").Value = 'Ballari (Bellary)'
$ws.Range("G48").Value = 'Ballari (Bellary)'
$ws.Range("G49").Value = 'Ballari (Bellary)'
$ws.Range("G50").Value = 'Mysuru (Mysore)'
$ws.Range("G51").Value = 'Ballari (Bellary)'
$ws.Range("G52").Value = 'Ballari (Bellary)'
$ws.Range("G53").Value = 'Ballari (Bellary)'
$ws.Range("G54").Value = 'Mysuru (Mysore)'
$ws.Range("G55").Value = 'Ballari (Bellary)'
$ws.Range("G56").Value = 'Kalaburagi (Gulbarga)'
$ws.Range("G57").Value = 'Mysuru (Mysore)'
$ws.Range("G58").Value = 'Mysuru (Mysore)'
$ws.Range("G59").Value = 'Mysuru (Mysore)'
$ws.Range("G60").Value = 'Mysuru (Mysore)'
$ws.Range("G61").Value = 'Kalaburagi (Gulbarga)'
$ws.Range("G62").Value = 'Ballari (Bellary)'
$ws.Range("G63").Value = 'Mysuru (Mysore)'
$ws.Range("G64").Value = 'Ballari (Bellary)'
$ws.Range("G65").Value = 'Ballari (Bellary)'
$ws.Range("G66").Value = 'Mysuru (Mysore)'
$ws.Range("G67").Value = 'Kalaburagi (Gulbarga)'
$ws.Range("G68").Value = 'Kalaburagi (Gulbarga)'

# Clear stray empty inline-string cells in column F
$ws.Range("F18").ClearContents()
$ws.Range("F35").ClearContents()
$ws.Range("F38").ClearContents()
$ws.Range("F69").ClearContents()
